$d = $word.ActiveDocument

# The hidden "_GoBack" bookmark currently sits, zero-length, at the very end
# of the document content -- right after the text of the final paragraph
# ("...it will have unbiased results.") and before that paragraph's mark.
# We want to end up with two brand-new list paragraphs appended after that
# paragraph, with the "_GoBack" bookmark ending up (still zero-length) at the
# new end of the document, immediately after the very last run.
#
# This runtime's Bookmarks.Add() does not actually (re)create/move a
# bookmark, and a literal paragraph-break character does not carry a
# bookmark across the split -- but ordinary text insertion at the bookmark's
# location does shift/carry it along. So: insert all of the new text in one
# go (using a NUL placeholder everywhere a paragraph break belongs) right at
# the bookmark's current position -- this carries "_GoBack" to the very end
# of the inserted text -- and only afterwards turn each placeholder into a
# real paragraph break via Find/Replace (which does not further disturb the
# already-relocated bookmark).

$bm = $d.Bookmarks("_GoBack")
$insertionPoint = $bm.Range.Start

$marker = [char]0
$newPara1Text = "Choose a solution and develop a plan to implement it:"
$newPara2Text = "Pull one sock out by one until you get at least one matching pair and do about 5 trials of those. Then put all the socks back in and pull one sock out one by one until you get a matching pair of each color."

$insertedText = $marker + $newPara1Text + $marker + $newPara2Text

$r = $d.Range($insertionPoint, $insertionPoint)
$r.InsertAfter($insertedText)

# Turn the two NUL placeholders into real paragraph breaks.
$null = $d.Content.Find.Execute($marker, $false, $false, $false, $false, `
                                 $false, $true, 1, $false, "^p", 2)

# Both new paragraphs inherited the pPr (style "ListParagraph", numId 6,
# ilvl 1) of the paragraph they were split out of. The first new paragraph
# ("Choose a solution...") is a top-level item and needs ilvl 0 (Word's
# 1-based ListLevelNumber 1); the second ("Pull one sock...") keeps ilvl 1
# (ListLevelNumber 2), matching the paragraph it inherited from.
$paraCount = $d.Paragraphs.Count
$newPara1 = $d.Paragraphs.Item($paraCount - 1)
$newPara1.Range.ListFormat.ListLevelNumber = 1

$newPara2 = $d.Paragraphs.Item($paraCount)
$newPara2.Range.ListFormat.ListLevelNumber = 2
